$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.809.30"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "1.758.34"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").Value = "'328.62"
$ws.Range("E5").Value = "  +1.40%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").Value = "'0.4515"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").Value = "'0.3492"
$ws.Range("E8").Value = "  -2.27%  "
$ws.Range("D9").Value = "'41.90"
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("D10").Value = "'0.07347"
$ws.Range("E10").Value = "  -1.75%  "
$ws.Range("D11").Value = "'1.090"
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").Value = "'20.57"
$ws.Range("E13").Value = "  -1.31%  "
$ws.Range("D14").Value = "'5.978"
$ws.Range("E14").Value = "  -0.36%  "
$ws.Range("D15").Value = "'7.167"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").Value = "1.765.23"
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").Value = "'91.87"
$ws.Range("E17").Value = "  -2.73%  "
$ws.Range("D18").Value = "'0.00001052"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "'0.06433"
$ws.Range("E19").Value = "  +0.71%  "
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("D21").Value = "'16.89"
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("D22").Value = "'5.764"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "27.894.77"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").Value = "'11.19"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").Value = "'2.159"
$ws.Range("E25").Value = "  +3.68%  "
$ws.Range("D26").Value = "'161.48"
$ws.Range("E26").Value = "  -2.50%  "
$ws.Range("D27").Value = "'20.09"
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("D28").Value = "1.969.95"
$ws.Range("E28").Value = "  +0.77%  "
$ws.Range("D29").Value = "'2.151"
$ws.Range("E29").Value = "  +0.69%  "
$ws.Range("D30").Value = "'123.48"
$ws.Range("E30").Value = "  -1.71%  "
$ws.Range("D31").Value = "'1.069"
$ws.Range("E31").Value = "  -1.55%  "
$ws.Range("D32").Value = "'0.09245"
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "'3.638"
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.545"
$ws.Range("E34").Value = "  +0.45%  "
$ws.Range("D35").Value = "'11.70"
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("D36").Value = "'0.02267"
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("D37").Value = "'0.06093"
$ws.Range("E37").Value = "  +1.50%  "
$ws.Range("D38").Value = "'0.2076"
$ws.Range("E38").Value = "  -0.68%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.6244"
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").Value = "'4.904"
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("D41").Value = "'1.181"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").Value = "'1.378"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("D43").Value = "'7.785"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("E44").Value = "  +0.71%  "
$ws.Range("D45").Value = "'3.730"
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("D46").Value = "'0.5804"
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("D47").Value = "'122.34"
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("D48").Value = "'1.923"
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("D49").Value = "'1.123"
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("D50").Value = "'0.06796"
$ws.Range("E50").Value = "  -1.35%  "
$ws.Range("D51").Value = "'72.57"
$ws.Range("E51").Value = "  +1.38%  "
